$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.2486567514670408
$ws.Range("E2").Value = -0.001477832512315258

$ws.Range("D3").Value = 0.4966880142591196
$ws.Range("E3").Value = 0.002891692954784419

$ws.Range("D4").Value = 0.09689601538188221
$ws.Range("E4").Value = 0.001949697796841399

$ws.Range("D5").Value = 0.1011801126389324
$ws.Range("E5").Value = 0.00697851509382641

$ws.Range("D6").Value = 0.05657910625302491
$ws.Range("E6").Value = 0.008552779653387521

$ws.Range("E7").Value = 0.002447709719560542
